# Ребалансировка Ставицкий.xlsx - "Add files via upload" edit
#
# The commit removes the workbook's external reference (the link to the
# external pricing workbook "[1]1") and bakes current VLOOKUP results into
# plain values, while also refreshing the "Текущая цена" (column G) market
# prices (and therefore the dependent "Текущая доходность" column N, which
# recalculates automatically via its shared formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# ---------------------------------------------------------------------
# 1) Refresh the "Текущая цена" (current price) column G with updated
#    market quotes.
# ---------------------------------------------------------------------
$newPrices = @{
    2  = 98.04
    3  = 99.95
    4  = 99.87
    5  = 96.49
    6  = 95.58
    7  = 97.19
    8  = 104.32
    9  = 96.06
    10 = 95.09
    11 = 99.8
    12 = 103.66
    13 = 109.82
    14 = 131.26
}

foreach ($row in $newPrices.Keys) {
    $ws.Cells.Item($row, 7).Value = $newPrices[$row]
}

# ---------------------------------------------------------------------
# 2) Capture the full-precision values currently produced by the
#    '[1]1'!VLOOKUP(...) formulas in column Q before the external link is
#    broken, so the baked-in constants keep their exact cached precision.
# ---------------------------------------------------------------------
$qValues = @{}
for ($row = 2; $row -le 13; $row++) {
    $qValues[$row] = $ws.Cells.Item($row, 17).Value2
}

# ---------------------------------------------------------------------
# 3) Break the external link to the pricing workbook. This removes the
#    <externalReferences> element from the workbook, deletes the cached
#    xl/externalLinks/externalLink1.xml part, and converts every formula
#    that referenced the external workbook into a plain cached value.
# ---------------------------------------------------------------------
$links = $wb.LinkSources(1)
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# ---------------------------------------------------------------------
# 4) Re-assert the exact pre-break values in column Q (BreakLink's own
#    conversion can lose precision), matching the previous cached results.
# ---------------------------------------------------------------------
foreach ($row in $qValues.Keys) {
    $ws.Cells.Item($row, 17).Value = $qValues[$row]
}

# ---------------------------------------------------------------------
# 5) Recalculate so the dependent formulas (e.g. N2:N10 = L/G*100) pick up
#    the new column G prices.
# ---------------------------------------------------------------------
$excel.Calculate()

# ---------------------------------------------------------------------
# 6) Cosmetic view-state changes captured in the same commit: row 2 is
#    shrunk back to its natural (63pt) height, and the frozen pane's
#    bottom-right selection is restored to B1.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 63

$ws.Range("B1").Select() | Out-Null
